# This script appends the new "OR / IN / BETWEEN / LIKE" note paragraphs to
# the end of the document, right before the final (empty) trailing paragraph,
# mirroring the target diff exactly (formatting, Wingdings arrow glyphs, the
# mid-paragraph page-break marker, etc. are all expressed as raw OOXML so the
# run/paragraph boundaries match the authored content precisely).

$d = $word.ActiveDocument

# The document always ends with one empty paragraph; the real last line of
# existing content is the paragraph right before it.
$lastIndex = $d.Paragraphs.Count
$anchorPara = $d.Paragraphs($lastIndex - 1)
$anchorText = $anchorPara.Range.Text
if ($anchorText -notmatch "Fetch the data which satisfied all the given conditions") {
    throw "Anchor paragraph not found where expected; got: $anchorText"
}

# Create a fresh paragraph right after the anchor (i.e. immediately before
# the trailing empty paragraph) to serve as the insertion point.
$insertRange = $anchorPara.Range
$insertRange.Collapse(0)
$insertRange.InsertParagraphAfter()

# The paragraph we just created now sits at the old trailing-empty-paragraph
# index; fill it (and grow it into many paragraphs) via InsertXML, which
# splices raw <w:p> OOXML directly into the body.
$targetPara = $d.Paragraphs($lastIndex)

$xmlFrag = '<w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>If we want fetch details where condition can satisfy either of values we should use OR.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Select * from emp where job = ‘MANAGER’ OR job = ‘SALESMAN’; </w:t></w:r><w:r><w:sym w:font="Wingdings" w:char="F0E0"/></w:r><w:r><w:t xml:space="preserve"> Fetch the records when job column has manager value or salesman value.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Exercise: Return names of</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> the employees where job is not manager nor salesman and sal &gt;=2000</w:t></w:r></w:p><w:p><w:r><w:t>Select ename from emp where job != ‘MANAGER’ AND job != ‘SALESMAN’ AND sal &gt;=2000;</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Exercise: Return names,hiring dates of the employees where location in Dallas or Chicago.</w:t></w:r></w:p><w:p><w:r><w:t>Select ename, hiredate from emp where deptno = 20 or deptno = 30;</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>IN Clause: To reduce the number of times to write deptno for OR Operator, we can use IN Clause</w:t></w:r></w:p><w:p><w:r><w:t>Select ename,h</w:t></w:r><w:r><w:t>iredate from emp where deptno</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> IN</w:t></w:r><w:r><w:t xml:space="preserve"> (20,30)</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Select ename, hiredate from emp where deptno </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>NOT IN</w:t></w:r><w:r><w:t xml:space="preserve"> (20,30) </w:t></w:r><w:r><w:sym w:font="Wingdings" w:char="F0E0"/></w:r><w:r><w:t xml:space="preserve"> Returns all the rows except the depno column values with 20,30 .</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>BETWEEN</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> OPERATOR: TO </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>filter</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> data in the given range we use BETWEEN.</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> we can use on numbers, dates and Textual data.</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>(INCLUSIVE)BETWEEN 1000 AND 2000 (1000 and 2000 included)</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Select * from EMP WHERE hiredate BETWEEN ‘05/01/1981’ AND ‘12/09/1982’; </w:t></w:r><w:r><w:sym w:font="Wingdings" w:char="F0E0"/></w:r><w:r><w:t xml:space="preserve"> Fetches records where hiredate is between the given range.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Select * from EMP WHERE sal NOT BETWEEN 600 AND 1500 </w:t></w:r><w:r><w:sym w:font="Wingdings" w:char="F0E0"/></w:r><w:r><w:t xml:space="preserve"> Fetches records which will not come under the given range condition for sal column.</w:t></w:r></w:p><w:p><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">select * from emp where comm is NULL; </w:t></w:r><w:r><w:sym w:font="Wingdings" w:char="F0E0"/></w:r><w:r><w:t xml:space="preserve"> Returns the records where comm column has null value.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Select * from emp where comm is NOT NULL </w:t></w:r><w:r><w:sym w:font="Wingdings" w:char="F0E0"/></w:r><w:r><w:t xml:space="preserve"> Returns the records where comm Column has values.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>USE Parenthesis () if we want to execute some conditions together.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">select * from emp </w:t></w:r><w:r><w:t xml:space="preserve">where </w:t></w:r><w:r><w:t>(</w:t></w:r><w:r><w:t>comm = 0</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">OR comm is NULL </w:t></w:r><w:r><w:t xml:space="preserve">) </w:t></w:r><w:r><w:t>AND sal BETWEEN 1101 AND 4999</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>AND sal &lt;&gt; 3000;</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">LIKE Operator: Wild Cards (%, *, </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>)</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr></w:p>'
$targetPara.Range.InsertXML($xmlFrag)

Write-Output ("Paragraphs after edit: " + $d.Paragraphs.Count)
